$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.8416365616161359, 1.1893813791480783, 1.3008920294914157, 2.907989791818589)"
$ws.Range("C2").Value = "NCT(2.4369115732160322, 1.7389440570237316, -0.010345543767976684, 4.325383411513661)"
$ws.Range("D2").Value = "NIG(1.6890916240832192, 1.2892561417156103, 0.8471986467003189, 3.4234820535667816)"
$ws.Range("E2").Value = "NCT(2.3136806570653907, 2.0460119801804924, -0.011397963140908587, 3.3272482753602803)"
